# Weekly update: insert two new rows of data at the top of the date-ordered
# block (rows 283/284), shifting the existing rows 283-329 down to 285-331.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the former row 283. Calling Insert()
# twice on the same row index pushes everything down by two rows total,
# copying formatting (e.g. the date style on column D) from the row below.
$ws.Rows.Item(283).Insert()
$ws.Rows.Item(283).Insert()

# New row 283 - "Primera" quality entry for the latest week (2023-01-25)
$ws.Cells.Item(283,1).Value2  = 11
$ws.Cells.Item(283,2).Value2  = 'Vega Monumental Concepción'
$ws.Cells.Item(283,3).Value2  = 'Bíobío'
$ws.Cells.Item(283,4).Value2  = 44951
$ws.Cells.Item(283,5).Value2  = 8
$ws.Cells.Item(283,6).Value2  = 100114013
$ws.Cells.Item(283,7).Value2  = 'Zanahoria'
$ws.Cells.Item(283,8).Value2  = 'Sin especificar'
$ws.Cells.Item(283,9).Value2  = 'Primera'
$ws.Cells.Item(283,10).Value2 = 800
$ws.Cells.Item(283,11).Value2 = 9000
$ws.Cells.Item(283,12).Value2 = 10000
$ws.Cells.Item(283,13).Value2 = 9500
$ws.Cells.Item(283,14).Value2 = '$/saco 20 kilos'
$ws.Cells.Item(283,15).Value2 = 'Región de Ñuble'
$ws.Cells.Item(283,16).Value2 = 475
$ws.Cells.Item(283,17).Value2 = 20
$ws.Cells.Item(283,18).Value2 = 'Hortaliza'

# New row 284 - "Segunda" quality entry for the same latest week
$ws.Cells.Item(284,1).Value2  = 11
$ws.Cells.Item(284,2).Value2  = 'Vega Monumental Concepción'
$ws.Cells.Item(284,3).Value2  = 'Bíobío'
$ws.Cells.Item(284,4).Value2  = 44951
$ws.Cells.Item(284,5).Value2  = 8
$ws.Cells.Item(284,6).Value2  = 100114013
$ws.Cells.Item(284,7).Value2  = 'Zanahoria'
$ws.Cells.Item(284,8).Value2  = 'Sin especificar'
$ws.Cells.Item(284,9).Value2  = 'Segunda'
$ws.Cells.Item(284,10).Value2 = 400
$ws.Cells.Item(284,11).Value2 = 8000
$ws.Cells.Item(284,12).Value2 = 8000
$ws.Cells.Item(284,13).Value2 = 8000
$ws.Cells.Item(284,14).Value2 = '$/saco 20 kilos'
$ws.Cells.Item(284,15).Value2 = 'Región de Ñuble'
$ws.Cells.Item(284,16).Value2 = 400
$ws.Cells.Item(284,17).Value2 = 20
$ws.Cells.Item(284,18).Value2 = 'Hortaliza'
